# "cierre 30 Dic 22"
#
# The workbook has two "vale" sheets (ARQUITECTO / VALES DE INSENTIVOS) plus
# an empty Hoja2. This closing edit bumps the ARQUITECTO vale's amount from
# $18,760.00 to $30,000.00 (and its spelled-out amount-in-words string to
# match), and leaves the sheet on ARQUITECTO as the active/selected tab
# (it had been left on "VALES DE INSENTIVOS").

$wb = $excel.ActiveWorkbook

$wsArquitecto = $wb.Worksheets.Item("ARQUITECTO        ")
$wsVales      = $wb.Worksheets.Item("VALES DE INSENTIVOS")

# New vale amount + matching amount-in-words text.
$wsArquitecto.Range("D1").Value = 30000
$wsArquitecto.Range("A2").Value = "   TREINTA   MIL      PESOS 00/100 M.N."

# Leave VALES DE INSENTIVOS' own selection untouched (still D6) but make
# ARQUITECTO the active sheet/tab with C17 selected, matching the closing
# state recorded in the workbook.
[void]$wsVales.Range("D6").Select()
[void]$wsArquitecto.Activate()
[void]$wsArquitecto.Range("C17").Select()
